$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": the two tracked files (d44ffe7d.. and e1a0821e..) swap
# rows 2/3, and the row that ends up holding the d44ffe7d file (now row 3)
# is updated to reflect it being ready for a fresh handoff.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "e1a0821e-2950-410f-ac19-156cb9e5b724.md"
$wsOverview.Range("B2").Value = "e2e\e1a0821e-2950-410f-ac19-156cb9e5b724.md"
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G2").Value = "2016-08-18 14:51:34"

$wsOverview.Range("A3").Value = "d44ffe7d-90fc-4235-9238-eb4b6785fa30.md"
$wsOverview.Range("B3").Value = "e2e\d44ffe7d-90fc-4235-9238-eb4b6785fa30.md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-18 14:52:41"

# Hyperlinks keep pointing at their original targets, but the display text
# needs to track the file name now shown in each row.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71c538d2a15f7f1b5599eb22f1aa86cecdda7bae/e2e/d44ffe7d-90fc-4235-9238-eb4b6785fa30.md", "", "", "e2e\e1a0821e-2950-410f-ac19-156cb9e5b724.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71c538d2a15f7f1b5599eb22f1aa86cecdda7bae/e2e/e1a0821e-2950-410f-ac19-156cb9e5b724.md", "", "", "e2e\d44ffe7d-90fc-4235-9238-eb4b6785fa30.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn": same row swap plus handoff-ready status for d44ffe7d.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "e1a0821e-2950-410f-ac19-156cb9e5b724.md"
$wsZh.Range("G2").Value = "e1a0821e-2950-410f-ac19-156cb9e5b724.17a387e0c195beaa3ef55d3ec3d993803b4a9f72.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-18 14:51:28"
$wsZh.Range("I2").Value = "e1a0821e-2950-410f-ac19-156cb9e5b724.md"
$wsZh.Range("J2").Value = "e1a0821e-2950-410f-ac19-156cb9e5b724.17a387e0c195beaa3ef55d3ec3d993803b4a9f72.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-18 14:51:55"

$wsZh.Range("A3").Value = "d44ffe7d-90fc-4235-9238-eb4b6785fa30.md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("G3").Value = "d44ffe7d-90fc-4235-9238-eb4b6785fa30.0341931d9c5303e347d81090f2db43f531e92132.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-18 14:52:36"
$wsZh.Range("I3").Value = "d44ffe7d-90fc-4235-9238-eb4b6785fa30.md"
$wsZh.Range("J3").Value = "d44ffe7d-90fc-4235-9238-eb4b6785fa30.0341931d9c5303e347d81090f2db43f531e92132.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-18 14:51:55"
$wsZh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71c538d2a15f7f1b5599eb22f1aa86cecdda7bae/e2e/d44ffe7d-90fc-4235-9238-eb4b6785fa30.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce2fdbb4928e577af105caca37c59dbb24865698/e2e/d44ffe7d-90fc-4235-9238-eb4b6785fa30.md."

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71c538d2a15f7f1b5599eb22f1aa86cecdda7bae/e2e/d44ffe7d-90fc-4235-9238-eb4b6785fa30.md", "", "", "e1a0821e-2950-410f-ac19-156cb9e5b724.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e29dc7d7e64e47e3646f0b54b2f6689acf6c6ce3/e2e/d44ffe7d-90fc-4235-9238-eb4b6785fa30.md", "", "", "e1a0821e-2950-410f-ac19-156cb9e5b724.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71c538d2a15f7f1b5599eb22f1aa86cecdda7bae/e2e/e1a0821e-2950-410f-ac19-156cb9e5b724.md", "", "", "d44ffe7d-90fc-4235-9238-eb4b6785fa30.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e29dc7d7e64e47e3646f0b54b2f6689acf6c6ce3/e2e/e1a0821e-2950-410f-ac19-156cb9e5b724.md", "", "", "d44ffe7d-90fc-4235-9238-eb4b6785fa30.md")

# Error Detail column needs more room to show the new message.
$wsZh.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# Sheet "de-de": same row swap plus handoff-ready status for d44ffe7d.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "e1a0821e-2950-410f-ac19-156cb9e5b724.md"
$wsDe.Range("G2").Value = "e1a0821e-2950-410f-ac19-156cb9e5b724.17a387e0c195beaa3ef55d3ec3d993803b4a9f72.de-de.xlf"
$wsDe.Range("I2").Value = "e1a0821e-2950-410f-ac19-156cb9e5b724.md"
$wsDe.Range("J2").Value = "e1a0821e-2950-410f-ac19-156cb9e5b724.17a387e0c195beaa3ef55d3ec3d993803b4a9f72.de-de.xlf"

$wsDe.Range("A3").Value = "d44ffe7d-90fc-4235-9238-eb4b6785fa30.md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("G3").Value = "d44ffe7d-90fc-4235-9238-eb4b6785fa30.0341931d9c5303e347d81090f2db43f531e92132.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-18 14:52:41"
$wsDe.Range("I3").Value = "d44ffe7d-90fc-4235-9238-eb4b6785fa30.md"
$wsDe.Range("J3").Value = "d44ffe7d-90fc-4235-9238-eb4b6785fa30.0341931d9c5303e347d81090f2db43f531e92132.de-de.xlf"
$wsDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71c538d2a15f7f1b5599eb22f1aa86cecdda7bae/e2e/d44ffe7d-90fc-4235-9238-eb4b6785fa30.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce2fdbb4928e577af105caca37c59dbb24865698/e2e/d44ffe7d-90fc-4235-9238-eb4b6785fa30.md."

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71c538d2a15f7f1b5599eb22f1aa86cecdda7bae/e2e/d44ffe7d-90fc-4235-9238-eb4b6785fa30.md", "", "", "e1a0821e-2950-410f-ac19-156cb9e5b724.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a2b693e1076937cf2ad8485ddf9b5ac63c6d2a25/e2e/d44ffe7d-90fc-4235-9238-eb4b6785fa30.md", "", "", "e1a0821e-2950-410f-ac19-156cb9e5b724.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71c538d2a15f7f1b5599eb22f1aa86cecdda7bae/e2e/e1a0821e-2950-410f-ac19-156cb9e5b724.md", "", "", "d44ffe7d-90fc-4235-9238-eb4b6785fa30.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a2b693e1076937cf2ad8485ddf9b5ac63c6d2a25/e2e/e1a0821e-2950-410f-ac19-156cb9e5b724.md", "", "", "d44ffe7d-90fc-4235-9238-eb4b6785fa30.md")

# Error Detail column needs more room to show the new message.
$wsDe.Columns.Item(16).ColumnWidth = 39.166666666666664
